$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing entry text in B2 (shared string index 1)
$ws.Range("B2").Value = "Le développement de l'interface graphique prend plus de temps que prévu, cependant la partie `"Chargeur de fichiers`" va finalement être abandonnée, car le processus est moins complexe que prévu. "

# Add a new row 3 with the date and new entry text
$ws.Range("A3").Value = 23.05
$ws.Range("B3").Value = "Après entretien avec M. Ithurbide, il a été décidé que la méthode de prendre un screenshot n'était effectivement pas optimisée, ou portable (celle-ci utilisait notamment un offset de coordonées hardodé, prévu pour ignorer spécifiquement les bordures de fenêtres windows 7). À la place, une conversion de l'image en matrice numpy sera utilisée"

# Update selection to match the diff
$ws.Range("B15").Select()
